$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.870.36'
$ws.Range('E2').Value = '  -1.67%  '
$ws.Range('D3').Value = '1.887.05'
$ws.Range('E3').Value = '  -2.62%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '0.7279'
$ws.Range('E5').Value = '  -5.73%  '
$ws.Range('D6').Value = '241.86'
$ws.Range('E6').Value = '  -1.83%  '
$ws.Range('D7').Value = '1.002'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '0.3084'
$ws.Range('E8').Value = '  -3.79%  '
$ws.Range('D9').Value = '26.27'
$ws.Range('E9').Value = '  -5.74%  '
$ws.Range('D10').Value = '0.06866'
$ws.Range('E10').Value = '  -2.52%  '
$ws.Range('D11').Value = '0.07939'
$ws.Range('E11').Value = '  -1.01%  '
$ws.Range('D12').Value = '0.7626'
$ws.Range('E12').Value = '  -2.40%  '
$ws.Range('D13').Value = '1.877.01'
$ws.Range('E13').Value = '  -3.10%  '
$ws.Range('D14').Value = '5.219'
$ws.Range('E14').Value = '  -2.67%  '
$ws.Range('D15').Value = '90.95'
$ws.Range('E15').Value = '  -3.96%  '
$ws.Range('D16').Value = '29.882.11'
$ws.Range('E16').Value = '  -1.65%  '
$ws.Range('D17').Value = '14.07'
$ws.Range('E17').Value = '  -2.55%  '
$ws.Range('D18').Value = '5.731'
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('D19').Value = '239.84'
$ws.Range('E19').Value = '  -6.28%  '
$ws.Range('D20').Value = '0.000007735'
$ws.Range('E20').Value = '  -2.54%  '
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').Value = '2.087.40'
$ws.Range('E22').Value = '  -4.77%  '
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').Value = '6.831'
$ws.Range('E24').Value = '  +1.53%  '
$ws.Range('D25').Value = '9.255'
$ws.Range('E25').Value = '  -3.01%  '
$ws.Range('D26').Value = '164.25'
$ws.Range('E26').Value = '  +0.36%  '
$ws.Range('D27').Value = '18.87'
$ws.Range('E27').Value = '  -1.09%  '
$ws.Range('D28').Value = '0.1263'
$ws.Range('E28').Value = '  -6.69%  '
$ws.Range('D29').Value = '2.005'
$ws.Range('E29').Value = '  -11.77%  '
$ws.Range('D30').Value = '1.346'
$ws.Range('E30').Value = '  -1.73%  '
$ws.Range('E31').Value = '  +0.66%  '
$ws.Range('D32').Value = '4.284'
$ws.Range('E32').Value = '  -2.89%  '
$ws.Range('D33').Value = '4.057'
$ws.Range('E33').Value = '  -1.55%  '
$ws.Range('D34').Value = '0.05057'
$ws.Range('E34').Value = '  -1.94%  '
$ws.Range('D35').Value = '1.266'
$ws.Range('E35').Value = '  -1.11%  '
$ws.Range('D36').Value = '0.7305'
$ws.Range('E36').Value = '  -2.30%  '
$ws.Range('D37').Value = '2.724'
$ws.Range('E37').Value = '  -2.13%  '
$ws.Range('D38').Value = '0.01911'
$ws.Range('E38').Value = '  -2.28%  '
$ws.Range('D39').Value = '2.767'
$ws.Range('E39').Value = '  -1.67%  '
$ws.Range('D40').Value = '6.315'
$ws.Range('E40').Value = '  -1.45%  '
$ws.Range('D41').Value = '74.22'
$ws.Range('E41').Value = '  -5.39%  '
$ws.Range('D42').Value = '0.4423'
$ws.Range('E42').Value = '  -1.89%  '
$ws.Range('D43').Value = '1.922'
$ws.Range('E43').Value = '  -2.52%  '
$ws.Range('D44').Value = '1.003'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').Value = '0.8331'
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '7.633'
$ws.Range('E46').Value = '  +1.68%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = '100.65'
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('D48').Value = '9.759'
$ws.Range('E48').Value = '  -0.45%  '
$ws.Range('D49').Value = '37.27'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').Value = '1.998.71'
$ws.Range('E50').Value = '  -4.18%  '
$ws.Range('D51').Value = '942.18'
$ws.Range('E51').Value = '  -4.10%  '
